# Applies the "Code Commit for Commodities" change:
# - adds two new rows (26, 27) to the Automation Tests sheet describing
#   AddCommodity_TC001 and EditCommodity_TC002 test cases
# - swaps the Run Mode (YES/No) values on rows 20 and 25
# - formats the new rows consistent with the rest of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap YES/No run-mode values on existing rows 20 and 25 ---
$ws.Range("C20").Value = "YES"
$ws.Range("C25").Value = "No"

# --- New row 26: AddCommodity_TC001 ---
$ws.Range("A26").Value = 'AddCommodity_TC001'
$ws.Range("B26").Value = 'Validate whehter Global admin is able to add new Shipper on following conditions.
a) Launch and login applcation as Global admin
b) Goto commodities and Click on add new commodity button.
C) Set Commodity name, upper limit, lower limit and pounds. 
d) click on save.
e) Customize Commodity webtable with required columns.
f) Search for record and check if details dispalyed corrrect.
g) Select record and click on delete button.
h) check whether record deleted successfully.'
$ws.Range("C26").Value = 'NO'
$ws.Range("D26").Value = 'Commodity added and deleted successfully'

# --- New row 27: EditCommodity_TC002 ---
$ws.Range("A27").Value = 'EditCommodity_TC002'
$ws.Range("B27").Value = 'Validate whehter Global admin is able to add new Shipper on following conditions.
a) Launch and login applcation as Global admin
b) Goto commodities and Click on add new commodity button.
C) Set Commodity name, upper limit, lower limit and pounds. 
d) click on save.
e) Customize Commodity webtable with required columns.
f) Search for record and check if details dispalyed corrrect.
g) Select record and edit with valid details.
h) check whether record edited successfully.'
$ws.Range("C27").Value = 'Yes'
$ws.Range("D27").Value = 'Commodity added and edited successfully'


# --- Formatting: column B wraps, columns A/C/D vertical-center ---
$ws.Range("A26").VerticalAlignment = -4108
$ws.Range("C26").VerticalAlignment = -4108
$ws.Range("D26").VerticalAlignment = -4108
$ws.Range("B26").WrapText = $true

$ws.Range("A27").VerticalAlignment = -4108
$ws.Range("C27").VerticalAlignment = -4108
$ws.Range("D27").VerticalAlignment = -4108
$ws.Range("B27").WrapText = $true

# --- Row heights matching the authored rows ---
$ws.Rows.Item(26).RowHeight = 150
$ws.Rows.Item(27).RowHeight = 150

# --- Restore active selection/view state ---
[void]$ws.Range("D27").Select()
